$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Reshape the sheet list:
#    before: [ODI Batting, ODI Bowling]
#    after:  [Player Info, ODI Batting, ODI Bowling, ODI Batting Extra]
# ---------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($batting)
$playerInfo.Name = "Player Info"

$bowling = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $bowling)
$extra.Name = "ODI Batting Extra"

# Re-fetch stable references (indices shifted after the inserts above).
$batting = $wb.Worksheets.Item("ODI Batting")
$bowling = $wb.Worksheets.Item("ODI Bowling")
$playerInfo = $wb.Worksheets.Item("Player Info")
$extra = $wb.Worksheets.Item("ODI Batting Extra")

function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
    $range.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------------
# 2. "Player Info" sheet - brand new.
# ---------------------------------------------------------------------------
$playerInfo.Range("A1:D2").NumberFormat = "@"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
Set-HeaderStyle $playerInfo.Range("A1:D1")

$playerInfo.Range("A2").Value = "4747"
$playerInfo.Range("B2").Value = "Asif Ali"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

# ---------------------------------------------------------------------------
# 3. "ODI Batting" sheet - rename MATCH_CARD_LINK -> MATCH_CODE, and turn the
#    scorecard URLs into bare match codes. Also drop the stray empty
#    INNING_NUMBER cells on "did not bat" rows.
# ---------------------------------------------------------------------------
$battingLastRow = $batting.Cells.Item($batting.Rows.Count, 1).End(-4162).Row

$batting.Columns.Item(4).NumberFormat = "@"
$batting.Range("D1").Value = "MATCH_CODE"

for ($row = 2; $row -le $battingLastRow; $row++) {
    $cell = $batting.Cells.Item($row, 4)
    $cell.Value = ($cell.Text -replace '^.*MatchCode=', '')

    $innings = $batting.Cells.Item($row, 2)
    if ($innings.Text -eq "") {
        $innings.ClearContents()
    }
}

# ---------------------------------------------------------------------------
# 4. "ODI Bowling" sheet - same MATCH_CARD_LINK -> MATCH_CODE treatment.
# ---------------------------------------------------------------------------
$bowlingLastRow = $bowling.Cells.Item($bowling.Rows.Count, 1).End(-4162).Row

$bowling.Columns.Item(2).NumberFormat = "@"
$bowling.Range("B1").Value = "MATCH_CODE"

for ($row = 2; $row -le $bowlingLastRow; $row++) {
    $cell = $bowling.Cells.Item($row, 2)
    $cell.Value = ($cell.Text -replace '^.*MatchCode=', '')
}

# ---------------------------------------------------------------------------
# 5. "ODI Batting Extra" sheet - brand new batting-detail table keyed by
#    MATCH_CODE.
# ---------------------------------------------------------------------------
$extra.Range("A1:F1").NumberFormat = "@"
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"
Set-HeaderStyle $extra.Range("A1:F1")

$extraRows = @(
    @("4174","6","","",""),
    @("4176","5","","",""),
    @("4177","3","5","3","12.53%"),
    @("4178","","","",""),
    @("4194","6","","",""),
    @("4197","6","0","1","5.56%"),
    @("4200","7","0","1","2.71%"),
    @("4201","6","1","2","12.66%"),
    @("4204","7","3","0","15.35%"),
    @("4227","","","",""),
    @("4287","","","",""),
    @("4292","4","2","4","14.13%"),
    @("4294","","","",""),
    @("4297","6","1","1","5.00%"),
    @("4300","7","2","0","7.41%"),
    @("4308","","","",""),
    @("4319","7","0","0","1.88%"),
    @("4458","6","0","0","0.73%"),
    @("4459","7","2","0","5.86%"),
    @("4567","5","","","")
)

$row = 2
foreach ($r in $extraRows) {
    $extra.Range("A" + $row).NumberFormat = "@"
    $extra.Range("A" + $row).Value = $r[0]

    if ($r[1] -eq "") {
        $extra.Cells.Item($row, 2).ClearContents()
    } else {
        $extra.Cells.Item($row, 2).NumberFormat = "General"
        $extra.Cells.Item($row, 2).Value = [double]$r[1]
    }

    "C","D","E" | ForEach-Object { } # no-op placeholder (kept for clarity)

    if ($r[2] -eq "") { $extra.Range("C" + $row).ClearContents() } else { $extra.Range("C" + $row).NumberFormat = "@"; $extra.Range("C" + $row).Value = $r[2] }
    if ($r[3] -eq "") { $extra.Range("D" + $row).ClearContents() } else { $extra.Range("D" + $row).NumberFormat = "@"; $extra.Range("D" + $row).Value = $r[3] }
    if ($r[4] -eq "") { $extra.Range("E" + $row).ClearContents() } else { $extra.Range("E" + $row).NumberFormat = "@"; $extra.Range("E" + $row).Value = $r[4] }

    $extra.Range("F" + $row).NumberFormat = "@"
    $extra.Range("F" + $row).Value = "NO"

    $row++
}

Write-Host "edit complete"
